$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 9.06
$ws.Range("C3").Value = 0.2420452041110066
$ws.Range("D3").Value = 0.06782578315039454
$ws.Range("E3").Value = 0.1025827849439245

# Row 4 (std)
$ws.Range("B4").Value = 0.238683256575942
$ws.Range("C4").Value = 0.3038493923259878
$ws.Range("D4").Value = 0.2470126788991783
$ws.Range("E4").Value = 0.2288250599216753

# Row 5 (min)
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 0.05980141166644689
$ws.Range("D5").Value = 0.0007213470278319244
$ws.Range("E5").Value = 0.01557299677543472

# Row 6 (25%)
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 0.1154084678727883
$ws.Range("D6").Value = 0.00226752011983907
$ws.Range("E6").Value = 0.02989674877322614

# Row 7 (50%)
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 0.1421636548960235
$ws.Range("D7").Value = 0.003182519236178738
$ws.Range("E7").Value = 0.03784719749647844

# Row 8 (75%)
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 0.2445616572524659
$ws.Range("D8").Value = 0.009606642761390709
$ws.Range("E8").Value = 0.06133239147452948

# Row 9 (max)
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 2.0601394884809
$ws.Range("D9").Value = 1.188572419413231
$ws.Range("E9").Value = 1
